# Elixir - Day 2 deck: turn the "Introduction to Mix" slide into a
# "Homework" slide with links, and drop the now-unused "Applications"
# and "Demo" slides that used to follow it (their notes page goes with
# the "Demo" slide automatically).

$p = $ppt.ActivePresentation

# --- Slide 26: "Introduction to Mix" -> "Homework" -------------------
$s = $p.Slides.Item(26)

$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Homework"

$body = $s.Shapes.Item(2)
$tr = $body.TextFrame.TextRange
$tr.Text = "https://github.com/georgiyolovski/elixir-workshop/tree/main/day2/homework"
$tr.InsertAfter("`r")
$tr.InsertAfter("`rThe following module references might be helpful:")
$tr.InsertAfter("`rhttps://hexdocs.pm/elixir/GenServer.html")

$p1 = $body.TextFrame.TextRange.Paragraphs(1, 1)
$p1.ActionSettings.Item(1).Hyperlink.Address = "https://github.com/georgiyolovski/elixir-workshop/tree/main/day2/homework"

$p4 = $body.TextFrame.TextRange.Paragraphs(4, 1)
$p4.IndentLevel = 2
$p4.ActionSettings.Item(1).Hyperlink.Address = "https://hexdocs.pm/elixir/GenServer.html"

# --- Remove the old "Applications" (27) and "Demo" (28) slides -------
# Delete from the back so the other index doesn't shift underneath us.
$p.Slides.Item(28).Delete()
$p.Slides.Item(27).Delete()
